$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Durée licence" (column H) wording for a few rows ---
$ws.Range("H2").Value = "1 an minimum"
$ws.Range("H3").Value = "1 mois minimum"
$ws.Range("H7").Value = "1 an minimum"

# --- Fill in "Autres fonctionnalités" (column I) for Syncfusion row ---
$ws.Range("I7").Value = "Autres actions sur les fichiers"

# --- Add new column J "Support", copying formatting from neighboring cells ---
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Support"

$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").Value = "Bon"

$ws.Range("H3:H6").Copy()
$ws.Range("J3:J6").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = "Moyen"

# --- Add new row 8: E-iceblue Spire (copy full formatting from row 7, then overwrite values) ---
$ws.Range("A7:J7").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)

$ws.Range("A8").Value = "E-iceblue Spire"
$ws.Range("B8").Value = "Bonne"
$ws.Range("C8").Value = "`$1,899/developer and deployment location               puis `$1,424.25 à renouveller pour les mises à jour"
$ws.Range("D8").Value = "MS Office: Word, Excel, PP"
$ws.Range("E8").Value = "Très bonne si optimisé"
$ws.Range("F8").Value = "?"
$ws.Range("G8").Value = "Difficile"
$ws.Range("H8").Value = "1 an minimum"
$ws.Range("I8").Value = "Autres actions sur les fichiers"
$ws.Range("J8").Value = "Bon"

$ws.Rows.Item(8).RowHeight = 58.8

# --- Add new row 9: CutePDF (name only) ---
$ws.Range("A4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("B9:J9").PasteSpecial(-4122)
$ws.Range("A9").Value = "CutePDF"

$ws.Rows.Item(9).RowHeight = 33

# --- Row height tweak for row 7 ---
$ws.Rows.Item(7).RowHeight = 64.8

# --- Column width tweaks for H, I, and the new J column ---
$ws.Columns.Item(8).ColumnWidth = 13.4375
$ws.Columns.Item(9).ColumnWidth = 22.96875
$ws.Columns.Item(10).ColumnWidth = 13.28125

# --- Update selection / view state ---
$ws.Range("C13").Select()
$excel.ActiveWindow.Zoom = 100
